$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.833.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.857.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5084'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('E8').Value = '  -2.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07164'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8888'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.67'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07528'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.856.37'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '91.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.234'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008532'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9999'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.875.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.012'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.092.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.443'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.821'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '146.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.044'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.640'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.672'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09223'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05094'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.065'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7341'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.149'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.195'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02004'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.073'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5311'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '118.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.500'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.395'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1473'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4630'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9994'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.953'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.559'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.87'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.76%  '
